$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AB2").Formula = "=K2"
